# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-12-09 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-12-10 Tuesday", 2)

# Update the division-problem table. The table has 5 "data" rows
# (1, 5, 9, 13, 17), each followed by 3 blank rows, and 5 columns.
# Each (row, col) pair below maps old text -> new text for that cell.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "65÷2=32, 1" },
    @{ Row = 1;  Col = 2; Text = "99÷6=16, 3" },
    @{ Row = 1;  Col = 3; Text = "68÷6=11, 2" },
    @{ Row = 1;  Col = 4; Text = "35÷3=11, 2" },
    @{ Row = 1;  Col = 5; Text = "99÷9=11, 0" },

    @{ Row = 5;  Col = 1; Text = "30÷8=3, 6" },
    @{ Row = 5;  Col = 2; Text = "13÷7=1, 6" },
    @{ Row = 5;  Col = 3; Text = "45÷2=22, 1" },
    @{ Row = 5;  Col = 4; Text = "51÷3=17, 0" },
    @{ Row = 5;  Col = 5; Text = "10÷8=1, 2" },

    @{ Row = 9;  Col = 1; Text = "67÷8=8, 3" },
    @{ Row = 9;  Col = 2; Text = "85÷7=12, 1" },
    @{ Row = 9;  Col = 3; Text = "32÷7=4, 4" },
    @{ Row = 9;  Col = 4; Text = "91÷6=15, 1" },
    @{ Row = 9;  Col = 5; Text = "44÷5=8, 4" },

    @{ Row = 13; Col = 1; Text = "86÷9=9, 5" },
    @{ Row = 13; Col = 2; Text = "59÷9=6, 5" },
    @{ Row = 13; Col = 3; Text = "97÷6=16, 1" },
    @{ Row = 13; Col = 4; Text = "64÷8=8, 0" },
    @{ Row = 13; Col = 5; Text = "47÷7=6, 5" },

    @{ Row = 17; Col = 1; Text = "71÷7=10, 1" },
    @{ Row = 17; Col = 2; Text = "38÷8=4, 6" },
    @{ Row = 17; Col = 3; Text = "42÷4=10, 2" },
    @{ Row = 17; Col = 4; Text = "99÷9=11, 0" },
    @{ Row = 17; Col = 5; Text = "71÷4=17, 3" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}

Write-Host "Applied date + $($updates.Count) table cell updates"
